# Adding PIB and Population data to data prep
# - Adds a new "Pop" worksheet (mnemonico/descricao mapping for population data)
#   at the end of the workbook, mirroring the existing "PIB" sheet's layout.
# - Moves the "selected/active tab" focus from "PIB" to the new "Pop" sheet.

$wb = $excel.ActiveWorkbook

# --- Update selection on the existing PIB sheet (tab focus moves off of it) ---
$pib = $wb.Worksheets.Item("PIB")
$pib.Range("A9").Select()

# --- Add the new "Pop" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$pop = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$pop.Name = "Pop"

# Header row (same mnemonico/descricao pattern used by the "PIB" sheet)
$pop.Cells.Item(1, 1).Value = "mnemonico"
$pop.Cells.Item(1, 2).Value = "descricao"

# Fill column B (descricao) first, then column A (mnemonico) - matches the
# order the values were authored in, which drives the shared-string table
# insertion order.
$pop.Cells.Item(2, 2).Value = "UF"
$pop.Cells.Item(3, 2).Value = "COD. UF"
$pop.Cells.Item(4, 2).Value = "COD. MUNIC"
$pop.Cells.Item(5, 2).Value = "NOME DO MUNICÍPIO"
$pop.Cells.Item(6, 2).Value = "POPULAÇÃO ESTIMADA"

$pop.Cells.Item(2, 1).Value = "UF"
$pop.Cells.Item(3, 1).Value = "CdUF"
$pop.Cells.Item(4, 1).Value = "CodIBGE"
$pop.Cells.Item(5, 1).Value = "Munip"
$pop.Cells.Item(6, 1).Value = "PopEstimada"

# Column A width, sized to fit its contents (mirrors the bestFit column
# widths used on the other sheets in this workbook).
$pop.Columns.Item(1).ColumnWidth = 19.6

# The newly-added sheet becomes the active / selected tab.
$pop.Activate()
$pop.Range("A1").Select()
